$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7083876132965088
$ws.Range("B1").Value = 1.350097894668579
$ws.Range("C1").Value = 4.137699127197266
$ws.Range("D1").Value = 2.277101755142212
$ws.Range("E1").Value = 0.7334194183349609
